$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.070.55"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").Value = "3.065.06"
$ws.Range("E3").Value = "  -1.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.75"
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.06"
$ws.Range("E6").Value = "  -1.09%  "

$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").Value = "3.062.86"
$ws.Range("E8").Value = "  -1.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  +2.72%  "

$ws.Range("E10").Value = "  +0.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.17"
$ws.Range("E11").Value = "  -11.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.497"
$ws.Range("E12").Value = "  +8.26%  "

$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.72"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").Value = "3.562.02"
$ws.Range("E15").Value = "  +1.03%  "

$ws.Range("D16").Value = "64.143.89"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").Value = "3.069.17"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("E18").Value = "  +1.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.80"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.48"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.96"
$ws.Range("E21").Value = "  +1.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.684"
$ws.Range("E22").Value = "  +1.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.44"
$ws.Range("E23").Value = "  +8.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.59"
$ws.Range("E24").Value = "  +1.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.58"
$ws.Range("E25").Value = "  +1.47%  "

$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.82"
$ws.Range("E27").Value = "  +0.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.16"
$ws.Range("E28").Value = "  +0.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.05"
$ws.Range("E29").Value = "  -1.33%  "

$ws.Range("E30").Value = "  +0.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.35"
$ws.Range("E31").Value = "  -0.14%  "

$ws.Range("E32").Value = "  -0.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.46"
$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.82"
$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.27"
$ws.Range("E35").Value = "  +2.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.93"
$ws.Range("E36").Value = "  -1.12%  "

$ws.Range("E37").Value = "  +1.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "451.33"
$ws.Range("E38").Value = "  -3.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0817"
$ws.Range("E39").Value = "  -2.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  +4.14%  "

$ws.Range("D41").Value = "3.026.96"
$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.30"
$ws.Range("E42").Value = "  -0.15%  "

$ws.Range("E43").Value = "  -1.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.267"
$ws.Range("E44").Value = "  +2.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "27.87"
$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.24"
$ws.Range("E46").Value = "  +7.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.113"
$ws.Range("E48").Value = "  +1.87%  "

$ws.Range("D49").Value = "0.0₃0520"
$ws.Range("E49").Value = "  -0.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "118.52"
$ws.Range("E50").Value = "  +1.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.09"
$ws.Range("E51").Value = "  -0.28%  "
